# Swap the "Valor Mora" (F column) values between row 16 and row 22 on Hoja1.
# Before: F16 = 52000, F22 = 32933
# After:  F16 = 32933, F22 = 52000

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("F16").Value = 32933
$ws.Range("F22").Value = 52000
